$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign literal text without Excel auto-converting date-like
# strings ("2025-03-30") into date serials, and without leaving a
# lingering non-default cell style behind.
function Set-TextCell($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Helper: assign a numeric value from its exact decimal-string form so
# the PowerShell tokenizer never has to parse a bare scientific-notation
# literal (e.g. 8.54e-05), which this host does not support directly.
function Set-NumCell($cell, $numStr) {
    $cell.Value = [double]$numStr
}

# --- Row 1: header labels (values unchanged by this update) ---
Set-TextCell $ws.Range("A1") "Row"
Set-TextCell $ws.Range("B1") "Prognose"
Set-TextCell $ws.Range("C1") "surveys"
Set-TextCell $ws.Range("D1") "production"
Set-TextCell $ws.Range("E1") "orders"
Set-TextCell $ws.Range("F1") "turnover"
Set-TextCell $ws.Range("G1") "financial"
Set-TextCell $ws.Range("H1") "labor market"
Set-TextCell $ws.Range("I1") "prices"
Set-TextCell $ws.Range("J1") "national accounts"
Set-TextCell $ws.Range("K1") "Revision"

# --- Rows 2-11: refreshed nowcast revision figures; row 12 is new ---
Set-TextCell $ws.Range("A2") "2025-03-30"
Set-NumCell $ws.Range("B2") "0.28698251968432625"
Set-NumCell $ws.Range("C2") "0"
Set-NumCell $ws.Range("D2") "0"
Set-NumCell $ws.Range("E2") "0"
Set-NumCell $ws.Range("F2") "0"
Set-NumCell $ws.Range("G2") "0"
Set-NumCell $ws.Range("H2") "0"
Set-NumCell $ws.Range("I2") "0"
Set-NumCell $ws.Range("J2") "0"
Set-NumCell $ws.Range("K2") "0"

Set-TextCell $ws.Range("A3") "2025-04-15"
Set-NumCell $ws.Range("B3") "0.28487027843248786"
Set-NumCell $ws.Range("C3") "0"
Set-NumCell $ws.Range("D3") "-0.0044333428824694151"
Set-NumCell $ws.Range("E3") "-0.0010399691952516956"
Set-NumCell $ws.Range("F3") "-0.00027128291572680853"
Set-NumCell $ws.Range("G3") "0.0011917950856776345"
Set-NumCell $ws.Range("H3") "8.5421416729136192e-05"
Set-NumCell $ws.Range("I3") "0.0034511262076645103"
Set-NumCell $ws.Range("J3") "0"
Set-NumCell $ws.Range("K3") "-0.0010959889684617474"

Set-TextCell $ws.Range("A4") "2025-04-30"
Set-NumCell $ws.Range("B4") "0.28777303148396577"
Set-NumCell $ws.Range("C4") "0.0064282723151962672"
Set-NumCell $ws.Range("D4") "0"
Set-NumCell $ws.Range("E4") "-6.1088935323815301e-05"
Set-NumCell $ws.Range("F4") "6.5101179895907441e-06"
Set-NumCell $ws.Range("G4") "0"
Set-NumCell $ws.Range("H4") "-0.00013038780871126357"
Set-NumCell $ws.Range("I4") "-0.0033918824176551211"
Set-NumCell $ws.Range("J4") "-2.8491252093075452e-08"
Set-NumCell $ws.Range("K4") "5.1358271234369646e-05"

Set-TextCell $ws.Range("A5") "2025-05-15"
Set-NumCell $ws.Range("B5") "0.28985021166392105"
Set-NumCell $ws.Range("C5") "0.01288589476996024"
Set-NumCell $ws.Range("D5") "-0.0040078828246083512"
Set-NumCell $ws.Range("E5") "-0.0020226936317724025"
Set-NumCell $ws.Range("F5") "0.001916409962867861"
Set-NumCell $ws.Range("G5") "-0.0078980133735725225"
Set-NumCell $ws.Range("H5") "-0.00034479010584411308"
Set-NumCell $ws.Range("I5") "0.0013289655705242779"
Set-NumCell $ws.Range("J5") "0"
Set-NumCell $ws.Range("K5") "0.00021928981240026957"

Set-TextCell $ws.Range("A6") "2025-05-30"
Set-NumCell $ws.Range("B6") "0.37999923310635941"
Set-NumCell $ws.Range("C6") "0.11689400803544725"
Set-NumCell $ws.Range("D6") "0"
Set-NumCell $ws.Range("E6") "-0.00015977753374392446"
Set-NumCell $ws.Range("F6") "0.0010414541543217765"
Set-NumCell $ws.Range("G6") "0"
Set-NumCell $ws.Range("H6") "0.00018464818064750261"
Set-NumCell $ws.Range("I6") "-0.027873517051278541"
Set-NumCell $ws.Range("J6") "0"
Set-NumCell $ws.Range("K6") "6.2205657044289708e-05"

Set-TextCell $ws.Range("A7") "2025-06-15"
Set-NumCell $ws.Range("B7") "0.37417052597053962"
Set-NumCell $ws.Range("C7") "0"
Set-NumCell $ws.Range("D7") "-0.014120110846221006"
Set-NumCell $ws.Range("E7") "-0.0031194278383382015"
Set-NumCell $ws.Range("F7") "-0.0025013153741022145"
Set-NumCell $ws.Range("G7") "0.012840295120547063"
Set-NumCell $ws.Range("H7") "0"
Set-NumCell $ws.Range("I7") "0.0018318590433654453"
Set-NumCell $ws.Range("J7") "0"
Set-NumCell $ws.Range("K7") "-0.00076000724107089246"

Set-TextCell $ws.Range("A8") "2025-06-30"
Set-NumCell $ws.Range("B8") "-0.050660471133130047"
Set-NumCell $ws.Range("C8") "-0.40409623933608529"
Set-NumCell $ws.Range("D8") "0"
Set-NumCell $ws.Range("E8") "2.684922383907323e-05"
Set-NumCell $ws.Range("F8") "-0.0015018551610716828"
Set-NumCell $ws.Range("G8") "0"
Set-NumCell $ws.Range("H8") "-0.00014024361475369252"
Set-NumCell $ws.Range("I8") "-0.023565134641280704"
Set-NumCell $ws.Range("J8") "0"
Set-NumCell $ws.Range("K8") "0.0044456264256826361"

Set-TextCell $ws.Range("A9") "2025-07-15"
Set-NumCell $ws.Range("B9") "-0.1668699988931493"
Set-NumCell $ws.Range("C9") "0"
Set-NumCell $ws.Range("D9") "-0.011858691367666251"
Set-NumCell $ws.Range("E9") "-0.003269092198539303"
Set-NumCell $ws.Range("F9") "-0.108583416822498"
Set-NumCell $ws.Range("G9") "0.0074945416448229446"
Set-NumCell $ws.Range("H9") "-0.0025429975446169674"
Set-NumCell $ws.Range("I9") "0.0026762982670196452"
Set-NumCell $ws.Range("J9") "0"
Set-NumCell $ws.Range("K9") "-0.00012616973854132074"

Set-TextCell $ws.Range("A10") "2025-07-30"
Set-NumCell $ws.Range("B10") "0.24430361815279766"
Set-NumCell $ws.Range("C10") "0.40435612756054184"
Set-NumCell $ws.Range("D10") "0"
Set-NumCell $ws.Range("E10") "-0.00049317307874291027"
Set-NumCell $ws.Range("F10") "0.0017683265799485017"
Set-NumCell $ws.Range("G10") "0"
Set-NumCell $ws.Range("H10") "-0.00049101839588120105"
Set-NumCell $ws.Range("I10") "0.034685653523882524"
Set-NumCell $ws.Range("J10") "-0.037866726417982277"
Set-NumCell $ws.Range("K10") "0.0092144272741804811"

Set-TextCell $ws.Range("A11") "2025-08-15"
Set-NumCell $ws.Range("B11") "0.49031312796468396"
Set-NumCell $ws.Range("C11") "0"
Set-NumCell $ws.Range("D11") "0.034420001535453479"
Set-NumCell $ws.Range("E11") "0.01911581259524308"
Set-NumCell $ws.Range("F11") "0.18877540207680757"
Set-NumCell $ws.Range("G11") "0.010193787644809046"
Set-NumCell $ws.Range("H11") "-0.00079631329278213045"
Set-NumCell $ws.Range("I11") "0.021183235549715407"
Set-NumCell $ws.Range("J11") "0"
Set-NumCell $ws.Range("K11") "-0.026882416297360157"

Set-TextCell $ws.Range("A12") "2025-08-30"
Set-NumCell $ws.Range("B12") "0.41286067996843312"
Set-NumCell $ws.Range("C12") "0.016895472620550807"
Set-NumCell $ws.Range("D12") "0"
Set-NumCell $ws.Range("E12") "0.00089726231537891766"
Set-NumCell $ws.Range("F12") "0.00017478520369325931"
Set-NumCell $ws.Range("G12") "0"
Set-NumCell $ws.Range("H12") "0.00012389032299918014"
Set-NumCell $ws.Range("I12") "-0.051686700233467157"
Set-NumCell $ws.Range("J12") "0"
Set-NumCell $ws.Range("K12") "-0.043857158225405868"

